# Updates the crypto price/volume snapshot (and re-ranks a few coins whose
# 1h change moved them past their neighbours), matching the Nov 1 2023
# GitHub Actions refresh of the cryptos list.
#
# Several "Price" cells (column D) hold values that look like plain numbers
# once the apostrophe thousands-separators fall away (e.g. "225.20",
# "0.0516"), so a bare .Value assignment would make Excel coerce them to
# numeric cells and drop the trailing zero / change their formatting. We
# force those through as text with a leading apostrophe and then restore
# the default "Normal" style so the cell itself carries no style index,
# matching the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.673.98"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.809.76"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'225.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'0.603"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'40.58"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.04%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("D12").Value = "2.070.19"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.808.72"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "'10.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "34.673.19"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'67.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "'241.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "0.0₃0769"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").Value = "'11.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'4.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("E24").Value = "  -3.22%  "
$ws.Range("D25").Value = "'172.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "'7.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("D27").Value = "'17.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").Value = "'0.0516"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.313.18"
$ws.Range("E35").Value = "  -4.37%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.642"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "'15.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.82%  "
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").Value = "'85.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.16%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0189"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  +5.86%  "
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "'0.944"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  +5.19%  "
$ws.Range("D47").Value = "1.969.10"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "'5.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'101.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "'0.0612"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
